$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Animal Count" is a new field being parsed out of the ARAMS data, placed
# right after "Herd Mark" (column O) and before "Animal Description"
# (old column P). Insert a new column there; this shifts the old P:R
# headers ("Animal Description", "Dept Country", "Dest Country") one
# column to the right, to Q:S, preserving their values/styles.
$ws.Columns.Item(16).Insert()

# New header cell for the inserted column.
$ws.Cells.Item(1, 16).Value = "Animal Count"

# The new header cell picks up the plain-default-font style (no explicit
# color) rather than the themed header style used by its neighbours -
# matching the formatting already used elsewhere in the sheet (the "AB123"
# data cell). Copy that formatting onto the new header cell.
$ws.Cells.Item(2, 2).Copy()
$ws.Cells.Item(1, 16).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
